$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.955.63"
$ws.Range("E2").Value = "  -0.89%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.547.44"
$ws.Range("E3").Value = "  +3.23%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.07"
$ws.Range("E5").Value = "  -0.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.71"
$ws.Range("E6").Value = "  +3.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("E8").Value = "  -1.91%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.545.35"
$ws.Range("E9").Value = "  +3.22%  "

$ws.Range("E10").Value = "  -1.41%  "

$ws.Range("E11").Value = "  -3.41%  "

$ws.Range("E12").Value = "  +0.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.352"
$ws.Range("E13").Value = "  -0.91%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.96"
$ws.Range("E14").Value = "  +1.83%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.000.71"
$ws.Range("E15").Value = "  +2.91%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.865.78"
$ws.Range("E16").Value = "  -0.90%  "

$ws.Range("E17").Value = "  -1.91%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.540.73"
$ws.Range("E18").Value = "  +2.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.46"
$ws.Range("E19").Value = "  +1.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "333.75"
$ws.Range("E20").Value = "  -2.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.26"
$ws.Range("E21").Value = "  -1.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.76"
$ws.Range("E22").Value = "  -1.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.66"
$ws.Range("E24").Value = "  -1.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  -4.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.59"
$ws.Range("E26").Value = "  +3.59%  "

$ws.Range("E28").Value = "  +11.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.32"
$ws.Range("E29").Value = "  +1.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.16"
$ws.Range("E30").Value = "  +4.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0808"
$ws.Range("E31").Value = "  -1.26%  "

$ws.Range("E32").Value = "  +0.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "177.04"
$ws.Range("E33").Value = "  +0.74%  "

$ws.Range("E34").Value = "  +3.51%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "406.62"
$ws.Range("E35").Value = "  +9.61%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.396"
$ws.Range("E36").Value = "  -1.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.78"
$ws.Range("E37").Value = "  -1.13%  "

$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.31"
$ws.Range("E39").Value = "  -4.10%  "

$ws.Range("E40").Value = "  +0.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.00"
$ws.Range("E42").Value = "  -3.33%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "151.19"
$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.73"
$ws.Range("E44").Value = "  -0.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.52"
$ws.Range("E45").Value = "  -1.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.602"
$ws.Range("E46").Value = "  +0.63%  "

$ws.Range("E47").Value = "  -0.88%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0517"
$ws.Range("E48").Value = "  -1.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0235"
$ws.Range("E49").Value = "  +3.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.36"
$ws.Range("E50").Value = "  +0.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.77"
$ws.Range("E51").Value = "  +1.15%  "
